# Applies the "New crime data collected" edits to CompStat_1 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text updates (report number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$c9 = $ws.Range("C9")
$c9.Characters(27,9).Text = "2/5/2024"
$c9.Characters(46,8).Text = "2/11/2024"

# --- Cells that were previously "n/a" (text) and are now numeric: set value + number format ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("K15").Value = 100
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D18").Value = 3
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = -33.333333333333
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 3
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 200
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M28").Value = 0
$ws.Range("M28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M29").Value = 0
$ws.Range("M29").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Cells that remain numeric: update value only (existing number format is retained) ---
$ws.Range("I15").Value = 2
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -33.333333333333
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = 43.75
$ws.Range("L16").Value = 9.523809523809
$ws.Range("M16").Value = -34.285714285714
$ws.Range("N16").Value = -80.833333333333
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 15.151515151515
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 17.021276595744
$ws.Range("M17").Value = 96.428571428571
$ws.Range("N17").Value = 37.5
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 5.555555555555
$ws.Range("L18").Value = -24
$ws.Range("M18").Value = -60.416666666666
$ws.Range("N18").Value = -87.074829931972
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -56.521739130434
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 7.692307692307
$ws.Range("I19").Value = 80
$ws.Range("J19").Value = 68
$ws.Range("K19").Value = 17.647058823529
$ws.Range("L19").Value = 42.857142857142
$ws.Range("M19").Value = 73.913043478260
$ws.Range("N19").Value = 29.032258064516
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 84
$ws.Range("L20").Value = 142.105263157895
$ws.Range("M20").Value = -6.122448979591
$ws.Range("N20").Value = -88.051948051948
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -17.777777777777
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = 16
$ws.Range("I21").Value = 226
$ws.Range("J21").Value = 178
$ws.Range("K21").Value = 26.966292134831
$ws.Range("L21").Value = 32.163742690058
$ws.Range("M21").Value = 8.653846153846
$ws.Range("N21").Value = -70.184696569920
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -2.061855670103
$ws.Range("I24").Value = 134
$ws.Range("J24").Value = 136
$ws.Range("K24").Value = -1.470588235294
$ws.Range("L24").Value = -8.219178082191
$ws.Range("M24").Value = 36.734693877551
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 5.263157894736
$ws.Range("I25").Value = 65
$ws.Range("J25").Value = 56
$ws.Range("K25").Value = 16.071428571428
$ws.Range("L25").Value = 16.071428571428
$ws.Range("M25").Value = -9.722222222222
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 40
